# Update the title text on the (only) slide of the InnerSource definition deck.
#   "VERSION 2 Inner "                 -> "VERSION 3 FOR PR "
#   "Source is the establishment of"   -> "Inner Source is the establishment of"
#
# The two phrases live in two separate runs inside the same paragraph of the
# "Rectangle 25" shape, so each is updated independently (matching the exact
# original text, including surrounding whitespace) to keep the run/formatting
# boundaries intact.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$tr.Replace("VERSION 2 Inner ", "VERSION 3 FOR PR ") | Out-Null
$tr.Replace("Source is the establishment of", "Inner Source is the establishment of") | Out-Null
